# Daily auto-push update (2026-02-25): two more time-of-day readings were
# recorded for "2026/02/25" (values 19 and 22) in addition to the existing
# 8/13/16 entries. They are appended to the end of the 2026/02/25 block,
# which sits right before the pre-existing "future" rows starting at
# 2026/12/29 - so those rows (and everything after them) shift down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first "2026/12/29" row (row 873),
# pushing all rows from 873 onward down to 875 onward.
$ws.Range("A873:D874").EntireRow.Insert()

# Seed the two new rows by copying the formatting/content of the row just
# above (row 872, the last "2026/02/25" / 16:00 entry) so the date & weekday
# columns stay plain text "2026/02/25" / "水" exactly like their neighbours,
# then overwrite just the time-of-day (column C) with the new readings.
$ws.Range("A872:D872").Copy()
$ws.Range("A873:D873").PasteSpecial()
$ws.Range("A872:D872").Copy()
$ws.Range("A874:D874").PasteSpecial()

$ws.Cells.Item(873, 3).Value = 19
$ws.Cells.Item(874, 3).Value = 22
